$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.098.46"
$ws.Range("E2").Value = "  +3.21%  "
$ws.Range("D3").Value = "1.656.88"
$ws.Range("E3").Value = "  +3.75%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.508"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +2.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0615"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0865"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").Value = "1.890.26"
$ws.Range("E12").Value = "  +3.75%  "
$ws.Range("D13").Value = "1.666.85"
$ws.Range("E13").Value = "  +4.22%  "
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("E15").Value = "  +3.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("D17").Value = "27.074.01"
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "238.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.89%  "
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("E22").Value = "  +5.02%  "
$ws.Range("E23").Value = "  +4.71%  "
$ws.Range("E24").Value = "  +3.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("E29").Value = "  +3.55%  "
$ws.Range("E31").Value = "  +1.96%  "
$ws.Range("D32").Value = "1.527.63"
$ws.Range("E32").Value = "  +4.02%  "
$ws.Range("E33").Value = "  +2.79%  "
$ws.Range("E34").Value = "  +4.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.14%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.579"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.892"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.67%  "
$ws.Range("E40").Value = "  +3.24%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.22%  "
$ws.Range("E43").Value = "  +3.64%  "
$ws.Range("D44").Value = "1.797.35"
$ws.Range("E44").Value = "  +3.58%  "
$ws.Range("E45").Value = "  +2.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.916"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("D48").Value = "0.0₆0105"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("E49").Value = "  +3.64%  "
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0978"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.38%  "
